# Germany Verbandsliga workbook update (14-04-2024)
# 1) Several existing rows had their home/away match order corrected -
#    for a handful of fixtures the row that used to carry match N's data
#    now carries match N+1's data and vice versa (id/date stay put, every
#    other column swaps).
# 2) One new fixture (row 137 / id 135) was appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $row1, $row2) {
    $cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $cell1 = $ws.Range($addr1)
        $cell2 = $ws.Range($addr2)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

# Row pairs whose fixture data needs to swap places
Swap-RowData $ws 3 4
Swap-RowData $ws 16 17
Swap-RowData $ws 46 47
Swap-RowData $ws 67 68
Swap-RowData $ws 69 70
Swap-RowData $ws 86 87
Swap-RowData $ws 125 126

# New fixture appended as row 137
$ws.Range("A137").Value2 = 135
$ws.Range("B137").Value2 = 8085441
$ws.Range("C137").Value2 = "Germany Verbandsliga"
$ws.Range("D137").Value2 = "Germany Verbandsliga"
$ws.Range("E137").Value2 = 45395.45833333334
$ws.Range("F137").Value2 = "TuS Hohenecken"
$ws.Range("G137").Value2 = "TB Jahn Zeiskam 1896"
$ws.Range("H137").Value2 = 1
$ws.Range("I137").Value2 = 4
$ws.Range("J137").Value2 = "A"
$ws.Range("K137").Value2 = 6
$ws.Range("L137").Value2 = 6
$ws.Range("M137").Value2 = 1.285
$ws.Range("N137").Value2 = 5.75
$ws.Range("O137").Value2 = 6
$ws.Range("P137").Value2 = 1.285
$ws.Range("Q137").Value2 = 1.75
$ws.Range("R137").Value2 = 1.85
$ws.Range("S137").Value2 = 1.95
$ws.Range("T137").Value2 = 3.5
$ws.Range("U137").Value2 = 2
$ws.Range("V137").Value2 = 1.8
$ws.Range("W137").Value2 = -1
$ws.Range("X137").Value2 = -1
$ws.Range("Y137").Value2 = 0.2849999999999999
$ws.Range("Z137").Value2 = -1
$ws.Range("AA137").Value2 = 0.95
$ws.Range("AB137").Value2 = 1
$ws.Range("AC137").Value2 = -1

# Match the formatting used by every other data row (bold/bordered id
# column, date-formatted date column)
$ws.Range("A136").Copy() | Out-Null
$ws.Range("A137").PasteSpecial(-4122) | Out-Null
$ws.Range("E136").Copy() | Out-Null
$ws.Range("E137").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
